$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Range("C3").ClearContents()
